$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 12.83529165041641
$ws.Cells.Item(2, 3).Value = 10.99820898269223
$ws.Cells.Item(2, 4).Value = 5.893046027424276
$ws.Cells.Item(2, 6).Value = 28.39980658462628
$ws.Cells.Item(2, 7).Value = 3.66181989566203
$ws.Cells.Item(2, 11).Value = 8.876563245456115
$ws.Cells.Item(2, 12).Value = 10.98729543644243
$ws.Cells.Item(2, 13).Value = 14.95039581915784
$ws.Cells.Item(2, 14).Value = 20.49586111256982
$ws.Cells.Item(2, 15).Value = 25.47995830748586

$ws.Cells.Item(3, 2).Value = 12.62795911426045
$ws.Cells.Item(3, 3).Value = 11.00615508401322
$ws.Cells.Item(3, 4).Value = 5.84811313882842
$ws.Cells.Item(3, 6).Value = 28.42654197996445
$ws.Cells.Item(3, 7).Value = 3.663582518899014
$ws.Cells.Item(3, 11).Value = 8.71739488862125
$ws.Cells.Item(3, 12).Value = 10.99485122418799
$ws.Cells.Item(3, 13).Value = 14.9222440815434
$ws.Cells.Item(3, 14).Value = 20.55648069228294
$ws.Cells.Item(3, 15).Value = 25.54184346480676

$ws.Cells.Item(4, 2).Value = 12.50156152051299
$ws.Cells.Item(4, 3).Value = 11.01146488230132
$ws.Cells.Item(4, 4).Value = 5.819841674439283
$ws.Cells.Item(4, 6).Value = 28.44940960452554
$ws.Cells.Item(4, 7).Value = 3.664723063385949
$ws.Cells.Item(4, 11).Value = 8.619750823744052
$ws.Cells.Item(4, 12).Value = 11.0011411446872
$ws.Cells.Item(4, 13).Value = 14.90731542632379
$ws.Cells.Item(4, 14).Value = 20.5954531107046
$ws.Cells.Item(4, 15).Value = 25.58457462372456

$ws.Cells.Item(5, 2).Value = 12.450348980082
$ws.Cells.Item(5, 3).Value = 11.01373738473826
$ws.Cells.Item(5, 4).Value = 5.808152944944826
$ws.Cells.Item(5, 6).Value = 28.46034968792664
$ws.Cells.Item(5, 7).Value = 3.665202546167804
$ws.Cells.Item(5, 11).Value = 8.580034816555314
$ws.Cells.Item(5, 12).Value = 11.00411994426996
$ws.Cells.Item(5, 13).Value = 14.90182929246391
$ws.Cells.Item(5, 14).Value = 20.61177643821007
$ws.Cells.Item(5, 15).Value = 25.6031767742559

$ws.Cells.Item(6, 2).Value = 12.44186500935834
$ws.Cells.Item(6, 3).Value = 11.01412130922019
$ws.Cells.Item(6, 4).Value = 5.806202010913964
$ws.Cells.Item(6, 6).Value = 28.4622641614983
$ws.Cells.Item(6, 7).Value = 3.665283053185416
$ws.Cells.Item(6, 11).Value = 8.573446084718022
$ws.Cells.Item(6, 12).Value = 11.00463968703284
$ws.Cells.Item(6, 13).Value = 14.90095453915497
$ws.Cells.Item(6, 14).Value = 20.61451363198805
$ws.Cells.Item(6, 15).Value = 25.60633740772278

$ws.Cells.Item(7, 2).Value = 12.50086956407002
$ws.Cells.Item(7, 3).Value = 11.01149508937614
$ws.Cells.Item(7, 4).Value = 5.819684710918913
$ws.Cells.Item(7, 6).Value = 28.44955058347197
$ws.Cells.Item(7, 7).Value = 3.664729470271689
$ws.Cells.Item(7, 11).Value = 8.619214826424896
$ws.Cells.Item(7, 12).Value = 11.00117963440681
$ws.Cells.Item(7, 13).Value = 14.90723901343711
$ws.Cells.Item(7, 14).Value = 20.59567146251066
$ws.Cells.Item(7, 15).Value = 25.58482068715545

$ws.Cells.Item(8, 2).Value = 12.76365608198011
$ws.Cells.Item(8, 3).Value = 11.00085959712187
$ws.Cells.Item(8, 4).Value = 5.877697043885567
$ws.Cells.Item(8, 6).Value = 28.40768538973194
$ws.Cells.Item(8, 7).Value = 3.662415577865996
$ws.Cells.Item(8, 11).Value = 8.821695777853419
$ws.Cells.Item(8, 12).Value = 10.98955843446499
$ws.Cells.Item(8, 13).Value = 14.94020299607633
$ws.Cells.Item(8, 14).Value = 20.51639992880107
$ws.Cells.Item(8, 15).Value = 25.50031308222642

$ws.Cells.Item(9, 2).Value = 13.28320702942388
$ws.Cells.Item(9, 3).Value = 10.98340571948064
$ws.Cells.Item(9, 4).Value = 5.985883682265392
$ws.Cells.Item(9, 6).Value = 28.37680367300908
$ws.Cells.Item(9, 7).Value = 3.658338468974407
$ws.Cells.Item(9, 11).Value = 9.217156770934967
$ws.Cells.Item(9, 12).Value = 10.97983957227181
$ws.Cells.Item(9, 13).Value = 15.02331555328651
$ws.Cells.Item(9, 14).Value = 20.37478958255547
$ws.Cells.Item(9, 15).Value = 25.37221026009399

$ws.Cells.Item(10, 2).Value = 13.66360994203471
$ws.Cells.Item(10, 3).Value = 10.97263445688907
$ws.Cells.Item(10, 4).Value = 6.061761214755014
$ws.Cells.Item(10, 6).Value = 28.38533022287313
$ws.Cells.Item(10, 7).Value = 3.655620835104691
$ws.Cells.Item(10, 11).Value = 9.503759844363694
$ws.Cells.Item(10, 12).Value = 10.98062655280868
$ws.Cells.Item(10, 13).Value = 15.09530823153653
$ws.Cells.Item(10, 14).Value = 20.27910290658143
$ws.Cells.Item(10, 15).Value = 25.30109941454344

$ws.Cells.Item(11, 2).Value = 13.83558869354145
$ws.Cells.Item(11, 3).Value = 10.96817527667377
$ws.Cells.Item(11, 4).Value = 6.095453619903076
$ws.Cells.Item(11, 6).Value = 28.39597107060909
$ws.Cells.Item(11, 7).Value = 3.65444422857579
$ws.Cells.Item(11, 11).Value = 9.632694616097789
$ws.Cells.Item(11, 12).Value = 10.98269510550026
$ws.Cells.Item(11, 13).Value = 15.1303542828751
$ws.Cells.Item(11, 14).Value = 20.23736886370936
$ws.Cells.Item(11, 15).Value = 25.27375790653196

$ws.Cells.Item(12, 2).Value = 13.90049588683952
$ws.Cells.Item(12, 3).Value = 10.96654970032439
$ws.Cells.Item(12, 4).Value = 6.108089874609425
$ws.Cells.Item(12, 6).Value = 28.40097015147334
$ws.Cells.Item(12, 7).Value = 3.654007210990479
$ws.Cells.Item(12, 11).Value = 9.681264963805575
$ws.Cells.Item(12, 12).Value = 10.98372325304396
$ws.Cells.Item(12, 13).Value = 15.14394820226314
$ws.Cells.Item(12, 14).Value = 20.22182197789605
$ws.Cells.Item(12, 15).Value = 25.26412520747629

$ws.Cells.Item(13, 2).Value = 13.88652761863512
$ws.Cells.Item(13, 3).Value = 10.96689700018179
$ws.Cells.Item(13, 4).Value = 6.105373930390054
$ws.Cells.Item(13, 6).Value = 28.39985043783039
$ws.Cells.Item(13, 7).Value = 3.654100951477651
$ws.Cells.Item(13, 11).Value = 9.670816510197739
$ws.Cells.Item(13, 12).Value = 10.98349095360541
$ws.Cells.Item(13, 13).Value = 15.1410062771416
$ws.Cells.Item(13, 14).Value = 20.22515887050544
$ws.Cells.Item(13, 15).Value = 25.26616770507851

$ws.Cells.Item(14, 2).Value = 13.84093336080884
$ws.Cells.Item(14, 3).Value = 10.96804027858483
$ws.Cells.Item(14, 4).Value = 6.096495684975631
$ws.Cells.Item(14, 6).Value = 28.39636293384613
$ws.Cells.Item(14, 7).Value = 3.654408104000358
$ws.Cells.Item(14, 11).Value = 9.636695888193
$ws.Cells.Item(14, 12).Value = 10.98277479200241
$ws.Cells.Item(14, 13).Value = 15.13146623861442
$ws.Cells.Item(14, 14).Value = 20.23608467001428
$ws.Cells.Item(14, 15).Value = 25.27295096534396

$ws.Cells.Item(15, 2).Value = 13.81297539835042
$ws.Cells.Item(15, 3).Value = 10.96874876591421
$ws.Cells.Item(15, 4).Value = 6.0910414628561
$ws.Cells.Item(15, 6).Value = 28.39435291549347
$ws.Cells.Item(15, 7).Value = 3.654597354446063
$ws.Cells.Item(15, 11).Value = 9.615761530343233
$ws.Cells.Item(15, 12).Value = 10.98236797066134
$ws.Cells.Item(15, 13).Value = 15.12566449233731
$ws.Cells.Item(15, 14).Value = 20.24281045910546
$ws.Cells.Item(15, 15).Value = 25.2771998194024

$ws.Cells.Item(16, 2).Value = 13.65234371359748
$ws.Cells.Item(16, 3).Value = 10.97293471271098
$ws.Cells.Item(16, 4).Value = 6.059542420850906
$ws.Cells.Item(16, 6).Value = 28.38477067079461
$ws.Cells.Item(16, 7).Value = 3.655698926075011
$ws.Cells.Item(16, 11).Value = 9.495300553119595
$ws.Cells.Item(16, 12).Value = 10.98052569364646
$ws.Cells.Item(16, 13).Value = 15.09306346689265
$ws.Cells.Item(16, 14).Value = 20.28186633373413
$ws.Cells.Item(16, 15).Value = 25.30298708539551

$ws.Cells.Item(17, 2).Value = 13.5534815138575
$ws.Cells.Item(17, 3).Value = 10.97561528044916
$ws.Cells.Item(17, 4).Value = 6.040004875917466
$ws.Cells.Item(17, 6).Value = 28.38062248853879
$ws.Cells.Item(17, 7).Value = 3.656389955283682
$ws.Cells.Item(17, 11).Value = 9.42099798480651
$ws.Cells.Item(17, 12).Value = 10.97983285365702
$ws.Cells.Item(17, 13).Value = 15.07364672727705
$ws.Cells.Item(17, 14).Value = 20.30628460614998
$ws.Cells.Item(17, 15).Value = 25.32008989846703

$ws.Cells.Item(18, 2).Value = 13.49652168237793
$ws.Cells.Item(18, 3).Value = 10.97719857611629
$ws.Cells.Item(18, 4).Value = 6.028689994603408
$ws.Cells.Item(18, 6).Value = 28.37887325626193
$ws.Cells.Item(18, 7).Value = 3.656793034804743
$ws.Cells.Item(18, 11).Value = 9.378128012237614
$ws.Cells.Item(18, 12).Value = 10.97959543156632
$ws.Cells.Item(18, 13).Value = 15.06269533367324
$ws.Cells.Item(18, 14).Value = 20.32049830591715
$ws.Cells.Item(18, 15).Value = 25.33039823928897

$ws.Cells.Item(19, 2).Value = 13.47722140439854
$ws.Cells.Item(19, 3).Value = 10.97774179206947
$ws.Cells.Item(19, 4).Value = 6.024845791568826
$ws.Cells.Item(19, 6).Value = 28.37839042015349
$ws.Cells.Item(19, 7).Value = 3.656930476724813
$ws.Cells.Item(19, 11).Value = 9.363591588550326
$ws.Cells.Item(19, 12).Value = 10.97954274326435
$ws.Cells.Item(19, 13).Value = 15.05902480244204
$ws.Cells.Item(19, 14).Value = 20.32533987324908
$ws.Cells.Item(19, 15).Value = 25.3339693742824

$ws.Cells.Item(20, 2).Value = 13.56401606395508
$ws.Cells.Item(20, 3).Value = 10.97532563651904
$ws.Cells.Item(20, 4).Value = 6.04209271911065
$ws.Cells.Item(20, 6).Value = 28.38099819184346
$ws.Cells.Item(20, 7).Value = 3.656315812922482
$ws.Cells.Item(20, 11).Value = 9.428921743837545
$ws.Cells.Item(20, 12).Value = 10.97988994427183
$ws.Cells.Item(20, 13).Value = 15.0756913077235
$ws.Cells.Item(20, 14).Value = 20.3036677613259
$ws.Cells.Item(20, 15).Value = 25.31822049451631

$ws.Cells.Item(21, 2).Value = 13.85433189553804
$ws.Cells.Item(21, 3).Value = 10.96770276253821
$ws.Cells.Item(21, 4).Value = 6.099106790248362
$ws.Cells.Item(21, 6).Value = 28.39736100881258
$ws.Cells.Item(21, 7).Value = 3.6543176544718
$ws.Cells.Item(21, 11).Value = 9.646725210504867
$ws.Cells.Item(21, 12).Value = 10.98297851107284
$ws.Cells.Item(21, 13).Value = 15.13425968133765
$ws.Cells.Item(21, 14).Value = 20.23286853485882
$ws.Cells.Item(21, 15).Value = 25.27093898496999

$ws.Cells.Item(22, 2).Value = 14.0427750913138
$ws.Cells.Item(22, 3).Value = 10.9630879288095
$ws.Cells.Item(22, 4).Value = 6.135654080784667
$ws.Cells.Item(22, 6).Value = 28.41370493485157
$ws.Cells.Item(22, 7).Value = 3.653061488268295
$ws.Cells.Item(22, 11).Value = 9.787569150856882
$ws.Cells.Item(22, 12).Value = 10.98642351890149
$ws.Cells.Item(22, 13).Value = 15.17441496247104
$ws.Cells.Item(22, 14).Value = 20.188093985975
$ws.Cells.Item(22, 15).Value = 25.24424025320274

$ws.Cells.Item(23, 2).Value = 13.94233827245193
$ws.Cells.Item(23, 3).Value = 10.96551747519228
$ws.Cells.Item(23, 4).Value = 6.116214715222499
$ws.Cells.Item(23, 6).Value = 28.4044660197495
$ws.Cells.Item(23, 7).Value = 3.653727389546428
$ws.Cells.Item(23, 11).Value = 9.712550537932781
$ws.Cells.Item(23, 12).Value = 10.98445474164196
$ws.Cells.Item(23, 13).Value = 15.15281410221203
$ws.Cells.Item(23, 14).Value = 20.21185441460043
$ws.Cells.Item(23, 15).Value = 25.25810508732794

$ws.Cells.Item(24, 2).Value = 13.55925377094883
$ws.Cells.Item(24, 3).Value = 10.97545645306254
$ws.Cells.Item(24, 4).Value = 6.041149061843375
$ws.Cells.Item(24, 6).Value = 28.38082635628668
$ws.Cells.Item(24, 7).Value = 3.656349314621687
$ws.Cells.Item(24, 11).Value = 9.425339883575379
$ws.Cells.Item(24, 12).Value = 10.97986363233396
$ws.Cells.Item(24, 13).Value = 15.07476629339062
$ws.Cells.Item(24, 14).Value = 20.30485029069323
$ws.Cells.Item(24, 15).Value = 25.31906417031005

$ws.Cells.Item(25, 2).Value = 13.14260935267284
$ws.Cells.Item(25, 3).Value = 10.98776549947546
$ws.Cells.Item(25, 4).Value = 5.95723492638852
$ws.Cells.Item(25, 6).Value = 28.37967286590186
$ws.Cells.Item(25, 7).Value = 3.65939244110431
$ws.Cells.Item(25, 11).Value = 9.110658388451332
$ws.Cells.Item(25, 12).Value = 10.98107326017896
$ws.Cells.Item(25, 13).Value = 14.9988864709818
$ws.Cells.Item(25, 14).Value = 20.41162582496985
$ws.Cells.Item(25, 15).Value = 25.40283010705001

